{"js": "// RBA v2.5 - Atualizacao da Tela\n// Replace the placeholder \"TRE\"/\"TERE\"/\"Tre\"/\"tre\" tokens with their\n// \"QWER\"/\"Qwer\"/\"Qewr\"/\"qwer\" counterparts (same casing pattern) in both\n// the document body (salutation line) and the primary page header\n// (address block).\n\n// 1) Body: \"A TERE, vem por meio desta convocar...\" -> \"A QWER, ...\"\nconst bodyResults = context.document.body.search(\"TERE\", {\n  matchCase: true,\n  matchWholeWord: false\n});\nawait context.sync();\nfor (let i = 0; i < bodyResults.items.length; i++) {\n  bodyResults.items[i].insertText(\"QWER\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Header: grab the primary header body of the (only) section.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst header = sections.items[0].getHeader(\"Primary\");\n\n// \"DIRETORIA DE ENSINO REGIAO TRE\" -> \"... QWER\"\nconst trePrimary = header.search(\"TRE\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nfor (let i = 0; i < trePrimary.items.length; i++) {\n  trePrimary.items[i].insertText(\"QWER\", \"Replace\");\n}\nawait context.sync();\n\n// \"TERE \u2013 DEP.\" -> \"QWER \u2013 DEP.\"\nconst tereHeader = header.search(\"TERE\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nfor (let i = 0; i < tereHeader.items.length; i++) {\n  tereHeader.items[i].insertText(\"QWER\", \"Replace\");\n}\nawait context.sync();\n\n// \"Tre, n\u00ba Tre \u2013 Tre \u2013 Tre \u2013 Tre\" -> \"Qwer, n\u00ba Qwer \u2013 Qewr \u2013 Qewr \u2013 Qwer\"\nconst treRuns = header.search(\"Tre\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nconst treReplacements = [\"Qwer\", \"Qwer\", \"Qewr\", \"Qewr\", \"Qwer\"];\nfor (let i = 0; i < treRuns.items.length; i++) {\n  treRuns.items[i].insertText(treReplacements[i], \"Replace\");\n}\nawait context.sync();\n\n// \"CEP: tre ... Tel: tre\" and \"Email: tre\" -> \"qwer\"\nconst treLower = header.search(\"tre\", { matchCase: true, matchWholeWord: false });\nawait context.sync();\nfor (let i = 0; i < treLower.items.length; i++) {\n  treLower.items[i].insertText(\"qwer\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# RBA v2.5 - Atualizacao da Tela\n# Replace the placeholder \"TRE\"/\"TERE\"/\"Tre\"/\"tre\" tokens with their\n# \"QWER\"/\"Qwer\"/\"Qewr\"/\"qwer\" counterparts (same casing pattern) in both\n# the document body (salutation line) and the primary page header\n# (address block).\n\n$wdReplaceOne      = [Microsoft.Office.Interop.Word.WdReplace]::wdReplaceOne\n$wdFindContinue    = [Microsoft.Office.Interop.Word.WdFindWrap]::wdFindContinue\n$wdHeaderFooterPrimary = [Microsoft.Office.Interop.Word.WdHeaderFooterIndex]::wdHeaderFooterPrimary\n\n$d = $word.ActiveDocument\n\n# 1) Body: \"A TERE, vem por meio desta convocar...\" -> \"A QWER, ...\"\n$bodyRange = $d.Content\n$bodyRange.Find.Execute(\"TERE\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"QWER\", $wdReplaceOne) | Out-Null\n\n# 2) Header: grab the primary header of the (only) section.\n$sec = $d.Sections.First\n$hdr = $sec.Headers.Item($wdHeaderFooterPrimary)\n\n# \"DIRETORIA DE ENSINO REGIAO TRE\" -> \"... QWER\"\n$hRange = $hdr.Range\n$hRange.Find.Execute(\"TRE\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"QWER\", $wdReplaceOne) | Out-Null\n\n# \"TERE - DEP.\" -> \"QWER - DEP.\"\n$hRange = $hdr.Range\n$hRange.Find.Execute(\"TERE\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"QWER\", $wdReplaceOne) | Out-Null\n\n# \"Tre, no Tre - Tre - Tre - Tre\" -> \"Qwer, no Qwer - Qewr - Qewr - Qwer\"\n$treReplacements = @(\"Qwer\", \"Qwer\", \"Qewr\", \"Qewr\", \"Qwer\")\n$hRange = $hdr.Range\nforeach ($rep in $treReplacements) {\n  $hRange.Find.Execute(\"Tre\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $rep, $wdReplaceOne) | Out-Null\n  $hRange.Collapse(0) | Out-Null\n  $hRange.End = $hdr.Range.End\n}\n\n# \"CEP: tre ... Tel: tre\" and \"Email: tre\" -> \"qwer\"\n$hRange = $hdr.Range\nfor ($i = 0; $i -lt 3; $i++) {\n  $hRange.Find.Execute(\"tre\", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, \"qwer\", $wdReplaceOne) | Out-Null\n  $hRange.Collapse(0) | Out-Null\n  $hRange.End = $hdr.Range.End\n}\n"}
